$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# Update item rows: Item_3/Item_4 -> puzzle block items
$ws.Range("B4").Value = "Item_PuzzleBlock_A"
$ws.Range("B5").Value = "Item_PuzzleBlock_B"

# Add new itemType description note in H2
$ws.Range("H2").Value = "itemType is an enum"

$ws.Range("C4").Value = "A puzzle block. Seems to be a part of some puzzle."
$ws.Range("C5").Value = "A puzzle block. Seems to be a part of some puzzle."

# Column width adjustments (mirrors the bestFit widths Excel computed for the
# widened B/C columns after the text changes above)
$ws.Columns.Item(2).ColumnWidth = 16.666666666666664
$ws.Columns.Item(3).ColumnWidth = 41.666666666666664

# Update selection to match authoring state
$ws.Range("D14").Select()
